# Refined metadata to be additional tab
#
# 1) Update the "time_taken" column (F) on the existing "data" sheet with
#    refreshed timestamps (the panel query was re-run later the same day).
# 2) Add a new "metadata" worksheet, placed after "data", carrying the
#    panel-query metadata (name/id/version/etc.) that used to only live
#    implicitly in the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "2021-10-05 14:20:34.181180"
$ws.Range("F3").Value = "2021-10-05 14:20:34.181188"
$ws.Range("F4").Value = "2021-10-05 14:20:34.181191"
$ws.Range("F5").Value = "2021-10-05 14:20:34.181194"
$ws.Range("F6").Value = "2021-10-05 14:20:34.181197"
$ws.Range("F7").Value = "2021-10-05 14:20:34.181199"
$ws.Range("F8").Value = "2021-10-05 14:20:34.181202"
$ws.Range("F9").Value = "2021-10-05 14:20:34.181204"
$ws.Range("F10").Value = "2021-10-05 14:20:34.181207"
$ws.Range("F11").Value = "2021-10-05 14:20:34.181210"
$ws.Range("F12").Value = "2021-10-05 14:20:34.181212"
$ws.Range("F13").Value = "2021-10-05 14:20:34.181215"
$ws.Range("F14").Value = "2021-10-05 14:20:34.181217"
$ws.Range("F15").Value = "2021-10-05 14:20:34.181219"
$ws.Range("F16").Value = "2021-10-05 14:20:34.181222"
$ws.Range("F17").Value = "2021-10-05 14:20:34.181224"
$ws.Range("F18").Value = "2021-10-05 14:20:34.181227"
$ws.Range("F19").Value = "2021-10-05 14:20:34.181229"
$ws.Range("F20").Value = "2021-10-05 14:20:34.181232"
$ws.Range("F21").Value = "2021-10-05 14:20:34.181234"
$ws.Range("F22").Value = "2021-10-05 14:20:34.181237"
$ws.Range("F23").Value = "2021-10-05 14:20:34.181239"
$ws.Range("F24").Value = "2021-10-05 14:20:34.181241"
$ws.Range("F25").Value = "2021-10-05 14:20:34.181244"
$ws.Range("F26").Value = "2021-10-05 14:20:34.181246"
$ws.Range("F27").Value = "2021-10-05 14:20:34.181249"
$ws.Range("F28").Value = "2021-10-05 14:20:34.181251"
$ws.Range("F29").Value = "2021-10-05 14:20:34.181254"
$ws.Range("F30").Value = "2021-10-05 14:20:34.181256"

# Add the new "metadata" sheet right after the "data" sheet.
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Glycogen storage disease"
$meta.Range("C2").Value = 528
# "1.7" is a panel *version* string, not a number -- force text storage
# (otherwise Excel auto-coerces "1.7" to the number 1.7) then drop back to
# the default "Normal" style so no stray style index sticks to the cell.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.7"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2021-08-26T10:55:47.645321Z"
$meta.Range("F2").Value = "2021-10-05 14:20:34.177465"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/528/?format=json"

# Match the bold/bordered header formatting already used on the "data"
# sheet's header row (style index 1) for the new header row + the leading
# index column, by copying formats across instead of guessing a style id.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Keep "data" as the active sheet/tab (adding "metadata" shouldn't change
# which tab is shown first).
$ws.Activate()
